$wb = $excel.ActiveWorkbook

# Template sheet with the row/column layout that matches the new market
# sheets (row heights, selection rectangle, border styling) -- Spain.
$template = $wb.Worksheets.Item("Spain")
$lastExisting = $wb.Worksheets.Item($wb.Worksheets.Count)

# New markets to add, in order, each with its Jira code (B4) and market label (B2)
$newMarkets = @(
    @{ Name = "Russia";  Code = "NGC-2929/T2925"; Market = "Russia Market"  },
    @{ Name = "Finland"; Code = "NGC-3130/T2887"; Market = "Finland Market" },
    @{ Name = "Hungary"; Code = "NGC-3104/T2979"; Market = "Hungary Market" }
)

$after = $lastExisting
foreach ($m in $newMarkets) {
    $template.Copy([System.Reflection.Missing]::Value, $after)
    $newSheet = $after.Next()
    $newSheet.Name = $m.Name
    $newSheet.Range("B4").Value = $m.Code
    $newSheet.Range("B2").Value = $m.Market
    $after = $newSheet
}

# The last added sheet becomes the active/selected tab
$after.Select()
